$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the D column (artefact name / author / reviewer) first
$ws.Range("D3").Value = "GeoBomber.py"
$ws.Range("D4").Value = "Éric Drouin"
$ws.Range("D5").Value = "Amar Hadjeres"

# Then the F column (project / dates)
$ws.Range("F3").Value = "GeoBomber"

# "03/10/2021" is ambiguous (parses as a valid M/D/Y date), so Excel would
# auto-convert a direct literal assignment into a date serial number.
# Enter it as a text formula first, then paste-special the value back in
# so it lands as plain text without picking up a numeric/date style.
$ws.Range("F4").Formula = "=""03/10/2021"""
$ws.Range("F4").Copy()
$ws.Range("F4").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("F5").Value = "29/09/2021"
$ws.Range("F6").Value = "17/10/2021"

# Préparation (min)
$ws.Range("D7").Value = 50

# Update the selected cell to match the saved view state
$ws.Range("F6").Select()
